# adj bpsk 1200 tune
$wb = $excel.ActiveWorkbook

$wsLoop = $wb.Worksheets.Item("LoopFilter LPF")
$wsBranch = $wb.Worksheets.Item("Branch LPF")
$wsGains = $wb.Worksheets.Item("Gains")

# LoopFilter LPF: scale factor 8 -> 1
$wsLoop.Range("B16").Value = 1

# Branch LPF: scale factor 1 -> 4
$wsBranch.Range("B18").Value = 4

# Update selections to match the author's saved cursor positions
$wsLoop.Range("B16").Select() | Out-Null
$wsBranch.Range("B28").Select() | Out-Null

# Branch LPF becomes the active/selected sheet (tab index 1)
$wsBranch.Activate()
